# presentation: mention cache locality, silly goose!
#
# Two edits to "P2P Apps.pptx":
#   1. Slide 26 ("PAST - Reclaim"), Content Placeholder: demote the
#      "Simply no longer guarantees it won't be deleted" bullet to the
#      second outline level.
#   2. Slide 38 ("Caching"), Content Placeholder: turn on Shrink-Text-
#      On-Overflow (normAutofit) and split the old "Cache Policy: ..."
#      bullet into three: a new "Cache locality" bullet, a new
#      sub-bullet explaining Pastry's proximity routing, and the
#      original "Cache Policy: GreeyDual-Size (GD-S)" bullet (minus the
#      "Cache " now living in the first bullet).

$p = $ppt.ActivePresentation

# --- Slide 26: "Simply ..." bullet becomes a sub-bullet (lvl 0 -> lvl 1) ---

$slide26 = $p.Slides.Item(26)
$contentShape26 = $slide26.Shapes.Item(3)
$tr26 = $contentShape26.TextFrame.TextRange

for ($i = 1; $i -le $tr26.Paragraphs().Count; $i++) {
    $para = $tr26.Paragraphs($i)
    if ($para.Text.StartsWith("Simply")) {
        # TextRange.IndentLevel is 1-based (1 == top level, matches
        # <a:pPr lvl="0"/> or no pPr at all); 2 == <a:pPr lvl="1"/>.
        $para.IndentLevel = 2
        break
    }
}

# --- Slide 38: rework the "Cache Policy:" bullet ---

$slide38 = $p.Slides.Item(38)
$contentShape38 = $slide38.Shapes.Item(3)

# Shrink text on overflow -> <a:bodyPr><a:normAutofit/></a:bodyPr>
$contentShape38.TextFrame.AutoSize = 2

$tr38 = $contentShape38.TextFrame.TextRange

$cacheParaIndex = -1
for ($i = 1; $i -le $tr38.Paragraphs().Count; $i++) {
    $para = $tr38.Paragraphs($i)
    if ($para.Text.StartsWith("Cache Policy:")) {
        $cacheParaIndex = $i
        break
    }
}

$cachePara = $tr38.Paragraphs($cacheParaIndex)

# Insert the two new bullets in front of the existing "Cache Policy: ..."
# paragraph. Splitting on a whole-paragraph range (rather than a
# Characters() sub-range) makes embedded `\r`s become real paragraph
# breaks instead of literal carriage returns.
[void]$cachePara.InsertBefore("Cache locality`rDue to Pastry" + [char]0x2019 + "s proximity`rCache ")

# The original run ("Cache Policy: ") is now glued onto the new "Cache "
# text, i.e. the paragraph reads "Cache Cache Policy: ...". Strip the
# duplicated leading "Cache " that belongs to the old run so it goes
# back to plain "Policy: ".
for ($i = 1; $i -le $tr38.Paragraphs().Count; $i++) {
    $para = $tr38.Paragraphs($i)
    if ($para.Text.StartsWith("Cache Cache Policy:")) {
        $dup = $tr38.Characters($para.Start + 6, 6)
        $dup.Text = ""
        break
    }
}

# New sub-bullet ("Due to Pastry's proximity") goes one level deeper.
for ($i = 1; $i -le $tr38.Paragraphs().Count; $i++) {
    $para = $tr38.Paragraphs($i)
    if ($para.Text.StartsWith("Due to Pastry")) {
        $para.IndentLevel = 2
        break
    }
}
